$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Living Standards Measurement Survey"
$ws.Range("B3").Value = "Computer Assisted Personal Interviewing"
$ws.Range("B7").Value = "Computer Assisted Telephone Interviewing"
$ws.Range("B13").Value = "Table of Contents"
$ws.Range("B15").Value = "Computer Assisted Web Interviewing"

$ws.Range("D14").Select()
